# Updated cryptos list (price + 1h volume change columns) per scheduled GitHub Actions refresh.
# Values in column D that are plain numeric strings are entered with a leading apostrophe so
# Excel keeps them as literal text (matching the source data's inlineStr cells) instead of
# auto-converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.858.06'
$ws.Range('E2').Value = '  +4.64%  '

$ws.Range('D3').Value = '3.274.82'
$ws.Range('E3').Value = '  +4.78%  '

$ws.Range('E4').Value = '  -0.03%  '

$ws.Range('D5').Value = '''581.41'
$ws.Range('E5').Value = '  +2.32%  '

$ws.Range('D6').Value = '''182.64'
$ws.Range('E6').Value = '  +8.45%  '

$ws.Range('E7').Value = '  -0.03%  '

$ws.Range('D8').Value = '''0.601'
$ws.Range('E8').Value = '  -0.10%  '

$ws.Range('D9').Value = '3.266.30'
$ws.Range('E9').Value = '  +4.42%  '

$ws.Range('E10').Value = '  +9.38%  '

$ws.Range('E11').Value = '  +3.73%  '

$ws.Range('D12').Value = '''0.419'
$ws.Range('E12').Value = '  +7.84%  '

$ws.Range('D13').Value = '3.837.97'
$ws.Range('E13').Value = '  +4.75%  '

$ws.Range('E14').Value = '  +1.22%  '

$ws.Range('D15').Value = '''28.67'
$ws.Range('E15').Value = '  +8.04%  '

$ws.Range('D16').Value = '67.813.91'
$ws.Range('E16').Value = '  +4.77%  '

$ws.Range('E17').Value = '  +5.43%  '

$ws.Range('D18').Value = '3.264.26'
$ws.Range('E18').Value = '  +4.32%  '

$ws.Range('E19').Value = '  +3.35%  '

$ws.Range('D20').Value = '''13.61'
$ws.Range('E20').Value = '  +7.49%  '

$ws.Range('D21').Value = '''376.36'
$ws.Range('E21').Value = '  +6.41%  '

$ws.Range('D22').Value = '''7.67'
$ws.Range('E22').Value = '  +6.64%  '

$ws.Range('E23').Value = '  +0.33%  '

$ws.Range('D24').Value = '''71.26'
$ws.Range('E24').Value = '  +3.67%  '

$ws.Range('E25').Value = '  +4.66%  '

$ws.Range('E26').Value = '  +6.40%  '

$ws.Range('D27').Value = '''9.66'
$ws.Range('E27').Value = '  +1.00%  '

$ws.Range('E28').Value = '  +3.69%  '

$ws.Range('E29').Value = '  -0.01%  '

$ws.Range('D30').Value = '''1.99'
$ws.Range('E30').Value = '  +4.54%  '

$ws.Range('E31').Value = '  +9.84%  '

$ws.Range('D32').Value = '''22.81'
$ws.Range('E32').Value = '  +5.21%  '

$ws.Range('B33').Value = 'Fetch.AI'
$ws.Range('C33').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D33').Value = '''1.28'
$ws.Range('E33').Value = '  +8.72%  '

$ws.Range('B34').Value = 'USDe'
$ws.Range('C34').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D34').Value = '''0.998'
$ws.Range('E34').Value = '  +0.01%  '

$ws.Range('D35').Value = '''6.95'
$ws.Range('E35').Value = '  +6.26%  '

$ws.Range('E36').Value = '  +6.62%  '

$ws.Range('D37').Value = '''163.28'
$ws.Range('E37').Value = '  +3.28%  '

$ws.Range('D38').Value = '''0.851'
$ws.Range('E38').Value = '  +3.18%  '

$ws.Range('E39').Value = '  +6.67%  '

$ws.Range('D40').Value = '''6.85'
$ws.Range('E40').Value = '  +12.92%  '

$ws.Range('B41').Value = 'EnergySwap'
$ws.Range('C41').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D41').Value = '''26.97'
$ws.Range('E41').Value = '  +3.16%  '

$ws.Range('B42').Value = 'Filecoin'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D42').Value = '''4.68'
$ws.Range('E42').Value = '  +13.21%  '

$ws.Range('E43').Value = '  +10.66%  '

$ws.Range('D44').Value = '''355.60'
$ws.Range('E44').Value = '  +12.60%  '

$ws.Range('D45').Value = '2.713.00'
$ws.Range('E45').Value = '  +2.96%  '

$ws.Range('D46').Value = '''25.55'
$ws.Range('E46').Value = '  +7.84%  '

$ws.Range('B47').Value = 'Hedera'
$ws.Range('C47').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D47').Value = '''0.0685'
$ws.Range('E47').Value = '  +5.74%  '

$ws.Range('B48').Value = 'OKB'
$ws.Range('C48').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D48').Value = '''40.93'
$ws.Range('E48').Value = '  +4.16%  '

$ws.Range('D49').Value = '''0.0282'
$ws.Range('E49').Value = '  +4.82%  '

$ws.Range('E50').Value = '  +8.14%  '

$ws.Range('E51').Value = '  +0.94%  '
